$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H99").Value = 3889.2856
$ws.Range("I99").Value = 175.66667
$ws.Range("K99").Value = 527.00001
$ws.Range("M99").Value = 970.99999
$ws.Range("H138").Value = 4172.9707
$ws.Range("J138").Value = 6668.25
$ws.Range("L138").Value = 20004.75
$ws.Range("N138").Value = -30284.75

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5639.6
$ws.Range("J2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("N2").Value = -4226
$ws.Range("H45").Value = 1558.6
$ws.Range("I45").Value = 1448.5
$ws.Range("J45").Value = 1999
$ws.Range("K45").Value = 1448.5
$ws.Range("L45").Value = 1999
$ws.Range("M45").Value = -1071.5
$ws.Range("N45").Value = -2753
$ws.Range("H63").Value = 2799.6667
$ws.Range("I63").Value = 2799.6667
$ws.Range("K63").Value = 2799.6667
$ws.Range("M63").Value = -2113.6667
$ws.Range("H66").Value = 2799.6667
$ws.Range("I66").Value = 2799.6667
$ws.Range("K66").Value = 13998.3335
$ws.Range("M66").Value = -10566.3335
$ws.Range("H74").Value = 1088938.8
$ws.Range("I74").Value = 1191790.1
$ws.Range("K74").Value = 1191790.1
$ws.Range("M74").Value = -1190916.1
$ws.Range("H77").Value = 1088938.8
$ws.Range("I77").Value = 1191790.1
$ws.Range("K77").Value = 5958950.5
$ws.Range("M77").Value = -5954582.5
$ws.Range("H80").Value = 184900
$ws.Range("J80").Value = 184900
$ws.Range("L80").Value = 184900
$ws.Range("N80").Value = -186896
$ws.Range("H83").Value = 184900
$ws.Range("J83").Value = 184900
$ws.Range("L83").Value = 554700
$ws.Range("N83").Value = -564684
$ws.Range("H97").Value = 1027.8206
$ws.Range("I97").Value = 946.8333
$ws.Range("K97").Value = 946.8333
$ws.Range("M97").Value = -450.8333
$ws.Range("H116").Value = 5639.6
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("N116").Value = -8588

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5639.6
$ws.Range("J3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("N3").Value = -4228
$ws.Range("H20").Value = 5730.7
$ws.Range("I20").Value = 7315.4614
$ws.Range("J20").Value = 2787.5715
$ws.Range("K20").Value = 7315.4614
$ws.Range("L20").Value = 2787.5715
$ws.Range("M20").Value = -7068.4614
$ws.Range("N20").Value = -3281.5715
$ws.Range("H86").Value = 37522.85
$ws.Range("I86").Value = 68870.5
$ws.Range("J86").Value = 6175.2
$ws.Range("K86").Value = 68870.5
$ws.Range("L86").Value = 6175.2
$ws.Range("M86").Value = -67747.5
$ws.Range("N86").Value = -8421.200000000001
$ws.Range("H89").Value = 37522.85
$ws.Range("I89").Value = 68870.5
$ws.Range("J89").Value = 6175.2
$ws.Range("K89").Value = 344352.5
$ws.Range("L89").Value = 30876
$ws.Range("M89").Value = -338736.5
$ws.Range("N89").Value = -42108
$ws.Range("H134").Value = 5265747
$ws.Range("I134").Value = 2513.0715
$ws.Range("K134").Value = 7539.2145
$ws.Range("M134").Value = -5004.2145

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31567798
$ws.Range("I31").Value = 35716380
$ws.Range("J31").Value = 2527727
$ws.Range("K31").Value = 35716380
$ws.Range("L31").Value = 2527727
$ws.Range("M31").Value = -35716085
$ws.Range("N31").Value = -2528317
$ws.Range("H34").Value = 31567798
$ws.Range("I34").Value = 35716380
$ws.Range("J34").Value = 2527727
$ws.Range("K34").Value = 35716380
$ws.Range("L34").Value = 2527727
$ws.Range("M34").Value = -35716178
$ws.Range("N34").Value = -2528131
$ws.Range("H75").Value = 72999.5
$ws.Range("J75").Value = 72999.5
$ws.Range("L75").Value = 72999.5
$ws.Range("N75").Value = -74995.5
$ws.Range("H78").Value = 72999.5
$ws.Range("J78").Value = 72999.5
$ws.Range("L78").Value = 218998.5
$ws.Range("N78").Value = -228982.5
$ws.Range("H94").Value = 1968.2727
$ws.Range("I94").Value = 2099.5
$ws.Range("K94").Value = 2099.5
$ws.Range("M94").Value = -1648.5
$ws.Range("H99").Value = 28176.25
$ws.Range("I99").Value = 13299.2
$ws.Range("K99").Value = 13299.2
$ws.Range("M99").Value = -11801.2
$ws.Range("H122").Value = 3821.2307
$ws.Range("I122").Value = 3671.7144
$ws.Range("K122").Value = 11015.1432
$ws.Range("M122").Value = -8565.143199999999
$ws.Range("H126").Value = 28176.25
$ws.Range("I126").Value = 13299.2
$ws.Range("K126").Value = 39897.60000000001
$ws.Range("M126").Value = -37427.60000000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4317.1787
$ws.Range("I55").Value = 2340.9092
$ws.Range("J55").Value = 5595.9414
$ws.Range("K55").Value = 7022.7276
$ws.Range("L55").Value = 16787.8242
$ws.Range("M55").Value = -6845.7276
$ws.Range("N55").Value = -17141.8242
$ws.Range("H68").Value = 2699.75
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 3333
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 9999
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -11621
$ws.Range("H71").Value = 2699.75
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 3333
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 29997
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -38109
$ws.Range("H105").Value = 21666.5
$ws.Range("I105").Value = 10000
$ws.Range("K105").Value = 30000
$ws.Range("M105").Value = -27379
$ws.Range("H107").Value = 5357784
$ws.Range("J107").Value = 7005211.5
$ws.Range("L107").Value = 21015634.5
$ws.Range("N107").Value = -21019474.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7684.115
$ws.Range("I70").Value = 6015.6665
$ws.Range("J70").Value = 9959.272000000001
$ws.Range("K70").Value = 6015.6665
$ws.Range("L70").Value = 9959.272000000001
$ws.Range("M70").Value = -5745.6665
$ws.Range("N70").Value = -10499.272
$ws.Range("H73").Value = 7684.115
$ws.Range("I73").Value = 6015.6665
$ws.Range("J73").Value = 9959.272000000001
$ws.Range("K73").Value = 6015.6665
$ws.Range("L73").Value = 9959.272000000001
$ws.Range("M73").Value = -5079.6665
$ws.Range("N73").Value = -11831.272
$ws.Range("H113").Value = 1326029.9
$ws.Range("J113").Value = 4632378
$ws.Range("L113").Value = 4632378
$ws.Range("N113").Value = -4636718
$ws.Range("H132").Value = 11114590
$ws.Range("J132").Value = 20003690
$ws.Range("L132").Value = 60011070
$ws.Range("N132").Value = -60016130

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4044.2
$ws.Range("J16").Value = 5399.6
$ws.Range("L16").Value = 5399.6
$ws.Range("N16").Value = -5739.6
$ws.Range("H22").Value = 12114
$ws.Range("I22").Value = 12872.637
$ws.Range("K22").Value = 12872.637
$ws.Range("M22").Value = -12577.637
$ws.Range("H27").Value = 12114
$ws.Range("I27").Value = 12872.637
$ws.Range("K27").Value = 12872.637
$ws.Range("M27").Value = -12765.637
$ws.Range("H40").Value = 4625
$ws.Range("H46").Value = 9999
$ws.Range("I46").Value = 9998
$ws.Range("K46").Value = 9998
$ws.Range("M46").Value = -9810
$ws.Range("H61").Value = 14751.25
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4298
$ws.Range("H93").Value = 1986767.9
$ws.Range("I93").Value = 656.6818
$ws.Range("J93").Value = 9269176
$ws.Range("K93").Value = 656.6818
$ws.Range("L93").Value = 9269176
$ws.Range("M93").Value = 591.3182
$ws.Range("N93").Value = -9271672
$ws.Range("H113").Value = 14751.25
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
$ws.Range("H132").Value = 5961.7
$ws.Range("I132").Value = 3230.5
$ws.Range("J132").Value = 7782.5
$ws.Range("K132").Value = 9691.5
$ws.Range("L132").Value = 23347.5
$ws.Range("M132").Value = -7161.5
$ws.Range("N132").Value = -28407.5
$ws.Range("H136").Value = 2810.5
$ws.Range("I136").Value = 2086.2222
$ws.Range("J136").Value = 4983.3335
$ws.Range("K136").Value = 6258.6666
$ws.Range("L136").Value = 14950.0005
$ws.Range("M136").Value = -3708.6666
$ws.Range("N136").Value = -20050.0005

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3699.0715
$ws.Range("I122").Value = 3284.8333
$ws.Range("J122").Value = 4009.75
$ws.Range("K122").Value = 9854.499899999999
$ws.Range("L122").Value = 12029.25
$ws.Range("M122").Value = -7404.499899999999
$ws.Range("N122").Value = -16929.25
$ws.Range("H132").Value = 359101.6
$ws.Range("I132").Value = 1973.5416
$ws.Range("J132").Value = 2501870
$ws.Range("K132").Value = 5920.6248
$ws.Range("L132").Value = 7505610
$ws.Range("M132").Value = -3390.6248
$ws.Range("N132").Value = -7510670
$ws.Range("H136").Value = 258638.72
$ws.Range("I136").Value = 2081.743
$ws.Range("J136").Value = 2503512.2
$ws.Range("K136").Value = 6245.228999999999
$ws.Range("L136").Value = 7510536.600000001
$ws.Range("M136").Value = -3695.228999999999
$ws.Range("N136").Value = -7515636.600000001
